$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.251.89"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "2.270.27"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "308.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  -1.16%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -1.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.11"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0806"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.13%  "
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("E13").Value = "  +0.90%  "
$ws.Range("D14").Value = "2.622.06"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "2.263.72"
$ws.Range("E16").Value = "  -1.63%  "
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "42.081.42"
$ws.Range("E18").Value = "  -1.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.71%  "
$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.58"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  -0.45%  "
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "164.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.25%  "
$ws.Range("E33").Value = "  +0.15%  "
$ws.Range("E34").Value = "  -0.29%  "
$ws.Range("E35").Value = "  -2.12%  "
$ws.Range("E36").Value = "  +0.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("E38").Value = "  -3.61%  "
$ws.Range("E39").Value = "  -3.78%  "
$ws.Range("E40").Value = "  -0.91%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.45%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.27"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -7.14%  "
$ws.Range("D43").Value = "1.959.49"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "18.91"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.74%  "
$ws.Range("E45").Value = "  -1.65%  "
$ws.Range("E46").Value = "  -3.24%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.84"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.14%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "53.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.91%  "
$ws.Range("D49").Value = "2.493.36"
$ws.Range("E49").Value = "  -0.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "92.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("E51").Value = "  -1.73%  "
